# Deploying to main from @ INTI-CMNB/kibot_variants_arduprog@cf435630e7e0ea3ad74ca84b569f64da0e434b7d
# - Update Revision / Date / KiCad Version fields on the BoM and DNF sheets.
# - Re-write the affected custom row heights so they serialize without a
#   trailing ".0" (e.g. "30.0" -> "30").

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("BoM", "DNF")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("D4").Value = '${git_hash}'
    $ws.Range("D5").Value = '${date}'
    $ws.Range("D6").Value = "6.0.11+dfsg-1~bpo11+1"
}

$bom = $wb.Worksheets.Item("BoM")
$bom.Rows.Item(13).RowHeight = 30
$bom.Rows.Item(15).RowHeight = 45
$bom.Rows.Item(16).RowHeight = 30
$bom.Rows.Item(17).RowHeight = 30

$dnf = $wb.Worksheets.Item("DNF")
$dnf.Rows.Item(11).RowHeight = 30
$dnf.Rows.Item(12).RowHeight = 45
$dnf.Rows.Item(13).RowHeight = 45
